# [MOSIP-14369] Fix: boolean values
#
# Column D (is_active) was previously stored as the formula TRUE() which
# evaluates to the numeric boolean 1. It should instead hold the literal
# text "TRUE" (a shared string), matching how the rest of the lookup
# columns (B, E, ...) are stored as text.
#
# We build one helper cell whose value is the text string "TRUE" (forcing
# a text/string result with a literal formula ="TRUE" that evaluates to a
# string, not a native boolean), then paste-special (values only) that
# string into every data row of column D (rows 2-132). Using copy +
# PasteSpecial(xlPasteValues) preserves each destination cell's existing
# style (s="1") instead of minting a brand new "quote-prefixed" style, and
# yields a genuine text cell (t="s") rather than Excel's usual TRUE/FALSE
# string -> boolean auto-coercion that a plain .Value assignment would
# trigger.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell, well outside the used range (A1:E132), used only to mint a
# text "TRUE" value we can copy from. Cleared again once we're done.
$helper = $ws.Cells.Item(1, 8)
$helper.Formula = '="TRUE"'
$helper.Copy()

$firstRow = 2
$lastRow = 132
for ($r = $firstRow; $r -le $lastRow; $r++) {
  $dst = $ws.Cells.Item($r, 4)
  $dst.PasteSpecial(-4163)
}

$helper.ClearContents()

# Restore the user's on-screen selection/scroll position recorded in the
# saved workbook: the view was scrolled down near the bottom of the sheet
# with the whole is_active column (minus the header) selected.
try {
  $win = $excel.ActiveWindow
  $win.ScrollRow = 109
  $win.ScrollColumn = 1
} catch {
}
$ws.Range("D2:D132").Select()
